# Updated cryptos list on Thu Feb 29 20:24:41 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "62.015.72"
$ws.Range("E2").Value = "  +2.72%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.408.16"
$ws.Range("E3").Value = "  +4.15%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("D5").Formula = "'406.55"
$ws.Range("E5").Value = "  +0.09%  "

# Row 6 - Solana
$ws.Range("D6").Formula = "'131.18"
$ws.Range("E6").Value = "  +18.41%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +8.26%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.04%  "

# Row 9 - Cardano
$ws.Range("D9").Formula = "'0.677"
$ws.Range("E9").Value = "  +10.13%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +13.00%  "

# Row 11 - Avalanche
$ws.Range("D11").Formula = "'42.23"
$ws.Range("E11").Value = "  +10.24%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.28%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.959.65"
$ws.Range("E13").Value = "  +4.27%  "

# Row 14 - Polkadot
$ws.Range("D14").Formula = "'8.57"
$ws.Range("E14").Value = "  +6.15%  "

# Row 15 - Chainlink
$ws.Range("D15").Formula = "'19.74"
$ws.Range("E15").Value = "  +4.95%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.403.13"
$ws.Range("E16").Value = "  +3.77%  "

# Row 17 - WrappedBTC (was Uniswap)
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "61.887.09"
$ws.Range("E17").Value = "  +2.50%  "

# Row 18 - Uniswap (was WrappedBTC)
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Formula = "'11.51"
$ws.Range("E18").Value = "  +9.97%  "

# Row 19 - Polygon
$ws.Range("E19").Value = "  +5.41%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  +19.09%  "

# Row 21 - ImmutableX
$ws.Range("E21").Value = "  +0.44%  "

# Row 22 - Litecoin
$ws.Range("D22").Formula = "'82.70"
$ws.Range("E22").Value = "  +13.51%  "

# Row 23 - InternetComputer(DFINITY)
$ws.Range("D23").Formula = "'13.18"
$ws.Range("E23").Value = "  +6.43%  "

# Row 24 - BitcoinCash
$ws.Range("D24").Formula = "'308.82"
$ws.Range("E24").Value = "  +4.76%  "

# Row 25 - PancakeSwap
$ws.Range("D25").Formula = "'3.17"
$ws.Range("E25").Value = "  +3.68%  "

# Row 26 - Filecoin
$ws.Range("E26").Value = "  +15.67%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Formula = "'29.82"
$ws.Range("E27").Value = "  +3.47%  "

# Row 28 - LEO
$ws.Range("D28").Formula = "'4.57"
$ws.Range("E28").Value = "  +7.31%  "

# Row 29 - RenderToken
$ws.Range("D29").Formula = "'7.49"
$ws.Range("E29").Value = "  +2.38%  "

# Row 30 - Kaspa
$ws.Range("D30").Formula = "'0.175"
$ws.Range("E30").Value = "  +1.95%  "

# Row 31 - Cosmos
$ws.Range("D31").Formula = "'11.87"
$ws.Range("E31").Value = "  +7.28%  "

# Row 32 - Hedera
$ws.Range("E32").Value = "  +3.84%  "

# Row 33 - Toncoin
$ws.Range("D33").Formula = "'2.62"
$ws.Range("E33").Value = "  +7.37%  "

# Row 34 - InjectiveProtocol
$ws.Range("D34").Formula = "'42.58"
$ws.Range("E34").Value = "  +10.79%  "

# Row 35 - Dai
$ws.Range("D35").Formula = "'1.00"
$ws.Range("E35").Value = "  -0.05%  "

# Row 36 - VeChain
$ws.Range("D36").Formula = "'0.0489"
$ws.Range("E36").Value = "  +2.76%  "

# Row 37 - OKB
$ws.Range("D37").Formula = "'52.37"
$ws.Range("E37").Value = "  +0.68%  "

# Row 38 - FirstDigitalUSD
$ws.Range("D38").Formula = "'0.997"
$ws.Range("E38").Value = "  -0.19%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  +5.38%  "

# Row 40 - Stacks
$ws.Range("E40").Value = "  -3.30%  "

# Row 41 - ARBITRUM
$ws.Range("D41").Formula = "'2.05"
$ws.Range("E41").Value = "  +9.98%  "

# Row 42 - Stellar
$ws.Range("E42").Value = "  +5.90%  "

# Row 43 - Monero
$ws.Range("D43").Formula = "'137.11"
$ws.Range("E43").Value = "  +1.98%  "

# Row 44 - NEARProtocol
$ws.Range("D44").Formula = "'3.99"
$ws.Range("E44").Value = "  +7.29%  "

# Row 45 - TheGraph (was Celestia)
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Formula = "'0.287"
$ws.Range("E45").Value = "  +0.05%  "

# Row 46 - Celestia (was TheGraph)
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").Formula = "'17.11"
$ws.Range("E46").Value = "  +6.23%  "

# Row 47 - WEMIXToken
$ws.Range("E47").Value = "  +1.76%  "

# Row 48 - EnergySwap
$ws.Range("D48").Formula = "'21.76"
$ws.Range("E48").Value = "  +4.85%  "

# Row 49 - Maker
$ws.Range("D49").Value = "2.152.48"
$ws.Range("E49").Value = "  +2.49%  "

# Row 50 - RocketPoolETH
$ws.Range("D50").Value = "3.744.23"
$ws.Range("E50").Value = "  +3.63%  "

# Row 51 - ApeXProtocol
$ws.Range("D51").Formula = "'2.34"
$ws.Range("E51").Value = "  -0.31%  "

Write-Output "done"
